$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.149.10"
$ws.Range("E2").Value = "  -0.92%  "
$ws.Range("D3").Value = "1.838.09"
$ws.Range("E3").Value = "  -1.45%  "
$ws.Range("D4").Value = "'0.9995"
$ws.Range("E4").Value = "  +0.18%  "
$ws.Range("D5").Value = "'240.53"
$ws.Range("E5").Value = "  -2.47%  "
$ws.Range("D6").Value = "'0.6851"
$ws.Range("E6").Value = "  -2.76%  "
$ws.Range("E7").Value = "  +0.21%  "
$ws.Range("E8").Value = "  -3.10%  "
$ws.Range("D9").Value = "'0.07412"
$ws.Range("E9").Value = "  -4.69%  "
$ws.Range("D10").Value = "'23.16"
$ws.Range("D11").Value = "'0.07643"
$ws.Range("D12").Value = "1.834.66"
$ws.Range("E12").Value = "  -1.47%  "
$ws.Range("D13").Value = "'5.054"
$ws.Range("E13").Value = "  -2.32%  "
$ws.Range("D14").Value = "'0.6804"
$ws.Range("E14").Value = "  -2.17%  "
$ws.Range("D15").Value = "'87.42"
$ws.Range("D16").Value = "'6.149"
$ws.Range("E16").Value = "  -7.47%  "
$ws.Range("D17").Value = "29.131.64"
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "'0.000008155"
$ws.Range("E18").Value = "  -2.67%  "
$ws.Range("D19").Value = "2.081.93"
$ws.Range("E19").Value = "  -0.60%  "
$ws.Range("D20").Value = "'229.18"
$ws.Range("E20").Value = "  -6.09%  "
$ws.Range("E21").Value = "  -2.36%  "
$ws.Range("D22").Value = "'0.9999"
$ws.Range("E22").Value = "  +0.13%  "
$ws.Range("D23").Value = "'7.346"
$ws.Range("E24").Value = "  +0.22%  "
$ws.Range("D25").Value = "'159.74"
$ws.Range("E25").Value = "  -0.21%  "
$ws.Range("E26").Value = "  -5.60%  "
$ws.Range("D27").Value = "'8.722"
$ws.Range("E27").Value = "  -2.45%  "
$ws.Range("E28").Value = "  -2.18%  "
$ws.Range("D29").Value = "'1.509"
$ws.Range("E29").Value = "  -2.11%  "
$ws.Range("D30").Value = "'4.262"
$ws.Range("E30").Value = "  +0.03%  "
$ws.Range("D31").Value = "'4.135"
$ws.Range("E31").Value = "  -2.10%  "
$ws.Range("D32").Value = "'1.194"
$ws.Range("E32").Value = "  -0.85%  "
$ws.Range("D33").Value = "'0.05252"
$ws.Range("E33").Value = "  +1.73%  "
$ws.Range("D34").Value = "'0.7548"
$ws.Range("E34").Value = "  -4.92%  "
$ws.Range("E35").Value = "  -4.30%  "
$ws.Range("E36").Value = "  -2.06%  "
$ws.Range("E37").Value = "  -0.28%  "
$ws.Range("D38").Value = "1.290.94"
$ws.Range("E38").Value = "  -3.33%  "
$ws.Range("D39").Value = "'0.01825"
$ws.Range("E39").Value = "  -2.81%  "
$ws.Range("D40").Value = "'2.719"
$ws.Range("E40").Value = "  -0.49%  "
$ws.Range("D41").Value = "'0.9360"
$ws.Range("E41").Value = "  -2.40%  "
$ws.Range("D42").Value = "'5.937"
$ws.Range("E42").Value = "  -2.30%  "
$ws.Range("D43").Value = "'104.63"
$ws.Range("E43").Value = "  -2.94%  "
$ws.Range("E44").Value = "  +0.08%  "
$ws.Range("B45").Value = "BabyDogeCoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D45").Value = "'0.00000000123"
$ws.Range("E45").Value = "  +4.61%  "
$ws.Range("B46").Value = "RocketPoolETH"
$ws.Range("C46").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D46").Value = "1.982.85"
$ws.Range("E46").Value = "  -0.69%  "
$ws.Range("B47").Value = "Mantle"
$ws.Range("C47").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D47").Value = "'0.5194"
$ws.Range("E47").Value = "  -0.02%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D48").Value = "'64.60"
$ws.Range("E48").Value = "  -1.15%  "
$ws.Range("D49").Value = "'9.461"
$ws.Range("E49").Value = "  -3.52%  "
$ws.Range("D50").Value = "'1.761"
$ws.Range("E50").Value = "  -1.14%  "
$ws.Range("D51").Value = "'0.07447"
$ws.Range("E51").Value = "  +17.49%  "
